$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-20 Sunday", "2025-07-21 Monday"),
    @("10÷3=", "62÷3="),
    @("14÷3=", "99÷3="),
    @("33÷7=", "24÷6="),
    @("18÷5=", "97÷5="),
    @("37÷9=", "21÷8="),
    @("70÷7=", "53÷8="),
    @("18÷2=", "32÷7="),
    @("55÷8=", "11÷6="),
    @("35÷5=", "79÷6="),
    @("92÷7=", "61÷3="),
    @("45÷5=", "36÷6="),
    @("24÷5=", "84÷9="),
    @("27÷6=", "32÷6="),
    @("84÷3=", "69÷6="),
    @("53÷2=", "13÷3="),
    @("11÷5=", "59÷5="),
    @("20÷5=", "10÷8="),
    @("17÷6=", "60÷6="),
    @("83÷9=", "64÷6="),
    @("90÷7=", "43÷8="),
    @("59÷2=", "81÷4="),
    @("47÷4=", "77÷8="),
    @("31÷2=", "70÷7="),
    @("54÷8=", "51÷6="),
    @("79÷3=", "84÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
